# Add a "3X" quantity column (K) to the BOM sheet: for each component row,
# compute the component count needed to build 3 sets (Quantity x 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell: K3 = "3X", formatted like the other header cells ---
$ws.Range("K3").Value = "3X"
$ws.Range("H3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

# --- Data rows 4..47: K = G * 3 ---
for ($r = 4; $r -le 47; $r++) {
    $ws.Range("K$r").Formula = "=G$r*3"
}

# Copy the Quantity column's formatting onto the new column so it matches
# (numeric style used throughout the sheet). Use a single source cell so
# every new-column cell gets the same uniform number style, same as a
# fill-down from K5 would produce.
$ws.Range("G5").Copy()
$ws.Range("K4:K47").PasteSpecial(-4122)

# Clear clipboard / marching ants
$excel.CutCopyMode = 0

# --- View state: selection moved to C44, no pinned top-left cell ---
$ws.Range("C44").Select()
